$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'309.86"
$ws.Range("E2").Value = "'0.12%"
$ws.Range("G2").Value = "'15"

# Row 3
$ws.Range("D3").Value = "'41.00"
$ws.Range("E3").Value = "'-0.80%"
$ws.Range("G3").Value = "'15"

# Row 4
$ws.Range("D4").Value = "'5.187"
$ws.Range("E4").Value = "'1.07%"
$ws.Range("G4").Value = "'15"

# Row 5
$ws.Range("E5").Value = "'0.40%"
$ws.Range("G5").Value = "'15"

# Row 6
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").Value = "'1.698"
$ws.Range("E6").Value = "'4.91%"
$ws.Range("G6").Value = "'15"

# Row 7
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "'0.9322"
$ws.Range("E7").Value = "'2.65%"
$ws.Range("G7").Value = "'15"

# Row 8
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.425"
$ws.Range("E8").Value = "'-3.14%"
$ws.Range("G8").Value = "'15"

# Row 9
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1281"
$ws.Range("E9").Value = "'12.12%"
$ws.Range("G9").Value = "'15"

# Row 10
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1838"
$ws.Range("E10").Value = "'1.67%"
$ws.Range("G10").Value = "'15"

# Row 11
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09139"
$ws.Range("E11").Value = "'0.10%"
$ws.Range("G11").Value = "'15"

# Row 12
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.04231"
$ws.Range("E12").Value = "'-0.30%"
$ws.Range("G12").Value = "'15"

# Row 13
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.1051"
$ws.Range("E13").Value = "'0.84%"
$ws.Range("G13").Value = "'15"

# Row 14
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001279"
$ws.Range("E14").Value = "'1.65%"
$ws.Range("G14").Value = "'15"

# Row 15
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005887"
$ws.Range("E15").Value = "'1.42%"
$ws.Range("G15").Value = "'15"

# Row 16
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.349"
$ws.Range("E16").Value = "'-0.17%"
$ws.Range("G16").Value = "'15"

# Row 17
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.299"
$ws.Range("E17").Value = "'0.71%"
$ws.Range("G17").Value = "'15"

# Row 18
$ws.Range("G18").Value = "'15"

# Row 19
$ws.Range("D19").Value = "'7.601"
$ws.Range("E19").Value = "'14.04%"
$ws.Range("G19").Value = "'15"

# Row 20
$ws.Range("D20").Value = "'0.1343"
$ws.Range("E20").Value = "'-1.62%"
$ws.Range("G20").Value = "'15"

# Row 21
$ws.Range("E21").Value = "'-0.56%"
$ws.Range("G21").Value = "'15"

# Row 22
$ws.Range("D22").Value = "'0.04018"
$ws.Range("E22").Value = "'-1.04%"
$ws.Range("G22").Value = "'15"

# Row 23
$ws.Range("D23").Value = "'0.001266"
$ws.Range("E23").Value = "'0.10%"
$ws.Range("G23").Value = "'15"

# Row 24
$ws.Range("D24").Value = "'0.004111"
$ws.Range("E24").Value = "'1.61%"
$ws.Range("G24").Value = "'15"

# Row 25
$ws.Range("E25").Value = "'0.11%"
$ws.Range("G25").Value = "'15"

# Row 26
$ws.Range("G26").Value = "'15"

# Row 27
$ws.Range("G27").Value = "'15"

# Row 28
$ws.Range("G28").Value = "'15"

# Row 29
$ws.Range("G29").Value = "'15"

# Row 30
$ws.Range("G30").Value = "'15"

# Row 31
$ws.Range("G31").Value = "'15"

# Row 32
$ws.Range("G32").Value = "'15"

# Row 33
$ws.Range("G33").Value = "'15"

# Row 34
$ws.Range("G34").Value = "'15"

# Row 35
$ws.Range("G35").Value = "'15"

# Row 36
$ws.Range("G36").Value = "'15"

# Row 37
$ws.Range("G37").Value = "'15"

# Row 38
$ws.Range("D38").Value = "'0.02550"
$ws.Range("E38").Value = "'4.72%"
$ws.Range("G38").Value = "'15"

# Row 39
$ws.Range("D39").Value = "'0.05312"
$ws.Range("E39").Value = "'0.86%"
$ws.Range("G39").Value = "'15"

# Row 40
$ws.Range("D40").Value = "'0.007843"
$ws.Range("E40").Value = "'0.67%"
$ws.Range("G40").Value = "'15"

# Row 41
$ws.Range("E41").Value = "'0.57%"
$ws.Range("G41").Value = "'15"

# Row 42
$ws.Range("D42").Value = "'0.006648"
$ws.Range("E42").Value = "'-1.06%"
$ws.Range("G42").Value = "'15"

# Row 43
$ws.Range("D43").Value = "'0.002202"
$ws.Range("E43").Value = "'13.55%"
$ws.Range("G43").Value = "'15"

# Row 44
$ws.Range("D44").Value = "'0.008119"
$ws.Range("E44").Value = "'7.58%"
$ws.Range("G44").Value = "'15"

# Row 45
$ws.Range("D45").Value = "'0.3087"
$ws.Range("E45").Value = "'0.23%"
$ws.Range("G45").Value = "'15"

# Row 46
$ws.Range("D46").Value = "'0.00006760"
$ws.Range("E46").Value = "'-2.08%"
$ws.Range("G46").Value = "'15"

# Row 47
$ws.Range("E47").Value = "'0.06%"
$ws.Range("G47").Value = "'15"

# Row 48
$ws.Range("D48").Value = "'0.2252"
$ws.Range("E48").Value = "'202.55%"
$ws.Range("G48").Value = "'15"

# Row 49
$ws.Range("E49").Value = "'3.43%"
$ws.Range("G49").Value = "'15"

# Row 50
$ws.Range("E50").Value = "'0.06%"
$ws.Range("G50").Value = "'15"

# Row 51
$ws.Range("E51").Value = "'0.06%"
$ws.Range("G51").Value = "'15"

